{"js": "// PK BTB Direksi Lajang - \"HASA MITRA\" -> \"HASAMITRA\" rebrand fix\n// (cetak sppk aman bulan 06 2024)\n//\n// Net textual changes made by the original edit:\n//   1. \"PT. BANK PERKREDITAN RAKYAT\" -> \"PT BANK PERKREDITAN RAKYAT\"  (drop the period after PT)\n//   2. \"HASA MITRA\"                  -> \"HASAMITRA\"                  (close up the company name, upper case)\n//   3. \"Hasa Mitra\"                  -> \"Hasamitra\"                  (close up the company name, title case)\n// The replaced fragments land in their own runs in the canonical OOXML, but the\n// Word object model only exposes text-level replace, so we do the text-level\n// equivalent here: a body-wide search & replace for each fragment.\n\nconst body = context.document.body;\n\nasync function replaceAll(searchText, replacement, options) {\n  const results = body.search(searchText, Object.assign({ matchCase: true }, options || {}));\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 1. \"PT. BANK PERKREDITAN RAKYAT\" -> \"PT BANK PERKREDITAN RAKYAT\" (2x)\nawait replaceAll(\"PT. BANK PERKREDITAN RAKYAT\", \"PT BANK PERKREDITAN RAKYAT\");\n\n// 2. \"HASA MITRA\" -> \"HASAMITRA\" (4x)\nawait replaceAll(\"HASA MITRA\", \"HASAMITRA\");\n\n// 3. \"Hasa Mitra\" -> \"Hasamitra\" (4x, inside the two \"setor langsung\" paragraphs)\nawait replaceAll(\"Hasa Mitra\", \"Hasamitra\");\n", "ps1": "# PK BTB Direksi Lajang - \"HASA MITRA\" -> \"HASAMITRA\" rebrand fix\n# (cetak sppk aman bulan 06 2024)\n#\n# Net textual changes made by the original edit:\n#   1. \"PT. BANK PERKREDITAN RAKYAT\" -> \"PT BANK PERKREDITAN RAKYAT\"  (drop the period after PT)\n#   2. \"HASA MITRA\"                  -> \"HASAMITRA\"                  (close up the company name, upper case)\n#   3. \"Hasa Mitra\"                  -> \"Hasamitra\"                  (close up the company name, title case)\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-AllText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n}\n\n# 1. \"PT. BANK PERKREDITAN RAKYAT\" -> \"PT BANK PERKREDITAN RAKYAT\" (2x)\nReplace-AllText \"PT. BANK PERKREDITAN RAKYAT\" \"PT BANK PERKREDITAN RAKYAT\"\n\n# 2. \"HASA MITRA\" -> \"HASAMITRA\" (4x)\nReplace-AllText \"HASA MITRA\" \"HASAMITRA\"\n\n# 3. \"Hasa Mitra\" -> \"Hasamitra\" (4x, inside the two \"setor langsung\" paragraphs)\nReplace-AllText \"Hasa Mitra\" \"Hasamitra\"\n"}
